$d = $word.ActiveDocument
$t = $d.Tables(1)

# Row 1: "100" -> "0M"
$t.Rows(1).Cells(1).Range.Text = "0M"

# Row 2: "0" -> "0M"
$t.Rows(2).Cells(1).Range.Text = "0M"

# Row 3: "161" -> "0M"
$t.Rows(3).Cells(1).Range.Text = "0M"

# Insert 10 new rows right after row 3, each with a single value
$newValues = @("102", "0.00002", "0.00006", "0.00004", "0.00001", "0.00004", "0.00004", "0.00005", "0.00426", "100.0")
$insertPos = 4
foreach ($val in $newValues) {
    $newRow = $t.Rows.Add($t.Rows($insertPos))
    $newRow.Cells(1).Range.Text = $val
    $insertPos = $insertPos + 1
}

# After inserting 10 rows, the former rows 34/35/36 are now 44/45/46.
# Row 44 (was a multi-run row "100<tab>0.00003<tab>...<tab>100.0") -> collapse to just "100"
$t.Rows(44).Cells(1).Range.Text = "100"

# Row 45 (was a multi-run row "2<tab>0.00002<tab>...<tab>100.0") -> collapse to just "0"
$t.Rows(45).Cells(1).Range.Text = "0"

# Row 46 (was an empty cell) -> set text to "161"
$t.Rows(46).Cells(1).Range.Text = "161"
